$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New columns F (custom_local_date_time_value_format) and G
# (custom_offset_date_time_value_format) with their header + sample rows.
# ---------------------------------------------------------------------------

# Headers (row 1) - same look & feel (font/format) as the existing E1 header.
$ws.Range("F1").Value = "custom_local_date_time_value_format"
$ws.Range("G1").Value = "custom_offset_date_time_value_format"
$ws.Range("E1:G1").NumberFormat = "General"

# Row 2
$ws.Range("F2").Value = "14-01-2022 13:00:12"
$ws.Range("F2").NumberFormat = "mm-dd-yy"
$ws.Range("G2").Value = "15 09 2025 15:40:37.180187 +02:00"
$ws.Range("G2").NumberFormat = "@"

# Row 3
$ws.Range("F3").Value = "15-01-2020 13:00:12"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("G3").Value = "15 09 2025 15:40:37.180187 +05:00"
$ws.Range("G3").NumberFormat = "@"

# Row 4 - A4/B4 now use the new custom LocalDateTime format.
$ws.Range("A4:B4").NumberFormat = "[$-F800]dddd, mmmm dd, yyyy"

# Row 5
$ws.Range("F5").NumberFormat = "m/d/yy h:mm;@"
$ws.Range("G5").Value = "15 09 2025 15:40:37.180187 Z"
$ws.Range("G5").NumberFormat = "General"

# A3/B3 now use the custom format applied to LocalDateTime values.
$ws.Range("A3:B3").NumberFormat = "mm-dd-yy"

# ---------------------------------------------------------------------------
# Sheet view / layout tweaks
# ---------------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 36.6640625
$ws.Columns.Item(7).ColumnWidth = 32.6640625

$ws.Range("F8").Select()

$ws.PageSetup.Orientation = 1
